# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (column D) values are plain text in this sheet (e.g. "70.151.69",
# "1.00"), not numbers, so numeric-looking updates are entered with a
# leading apostrophe to force Excel to keep them as text instead of
# silently parsing them into numeric values (which would, e.g., turn
# "1.00" into 1 or reformat "569.85000000000002").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.159.38"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "'3.530.99"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'569.85"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("D6").Value = "'182.90"
$ws.Range("E6").Value = "  -5.15%  "
$ws.Range("D7").Value = "'3.516.41"
$ws.Range("E7").Value = "  -2.49%  "
$ws.Range("D8").Value = "'0.613"
$ws.Range("E8").Value = "  -4.05%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'0.182"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "'0.640"
$ws.Range("E11").Value = "  -4.06%  "
$ws.Range("D12").Value = "'53.53"
$ws.Range("E12").Value = "  -6.48%  "
$ws.Range("D13").Value = "'0.0000297"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").Value = "'9.43"
$ws.Range("E14").Value = "  -4.19%  "
$ws.Range("D15").Value = "'4.119.40"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").Value = "'19.24"
$ws.Range("E16").Value = "  -5.58%  "
$ws.Range("D17").Value = "'3.560.48"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "'69.283.45"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "'12.24"
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").Value = "'1.03"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").Value = "'500.13"
$ws.Range("E22").Value = "  +3.01%  "
$ws.Range("D23").Value = "'19.40"
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("D24").Value = "'4.88"
$ws.Range("E24").Value = "  -4.36%  "
$ws.Range("D25").Value = "'4.29"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").Value = "'93.54"
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("D27").Value = "'11.23"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").Value = "'2.91"
$ws.Range("E28").Value = "  -6.78%  "
$ws.Range("D29").Value = "'9.13"
$ws.Range("E29").Value = "  -3.41%  "
$ws.Range("D30").Value = "'31.22"
$ws.Range("E30").Value = "  -4.18%  "
$ws.Range("D31").Value = "'7.51"
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("D32").Value = "'12.28"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'65.15"
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("D34").Value = "'0.114"
$ws.Range("E34").Value = "  -7.60%  "
$ws.Range("D35").Value = "'565.32"
$ws.Range("E35").Value = "  -7.01%  "
$ws.Range("D36").Value = "'3.09"
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'0.997"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'37.79"
$ws.Range("E38").Value = "  -5.75%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.398"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'3.27"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "'0.0₃0775"
$ws.Range("E41").Value = "  -6.54%  "
$ws.Range("D42").Value = "'3.36"
$ws.Range("E42").Value = "  -5.54%  "
$ws.Range("D43").Value = "'0.133"
$ws.Range("E43").Value = "  -9.79%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'3.56"
$ws.Range("E44").Value = "  +5.42%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.96"
$ws.Range("E45").Value = "  -5.38%  "
$ws.Range("D46").Value = "'0.0439"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("D47").Value = "'3.163.75"
$ws.Range("E47").Value = "  -4.75%  "
$ws.Range("D48").Value = "'9.23"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("D49").Value = "'0.134"
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.16"
$ws.Range("E51").Value = "  -1.74%  "
